# Regenerate merged AHB files
#
# 1. Rename the "_old" / "_new" header-column suffixes to the new
#    version-tag suffixes "_FV2404" / "_FV2410" (row 1, columns A:U).
# 2. Turn the A1:U57 range into a real Excel Table ("Table1").
# 3. Freeze the header row (pane split under row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header labels -------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$headerRange.Replace("_old", "_FV2404", 2, 1, $false, $false, $false, $false) | Out-Null
$headerRange.Replace("_new", "_FV2410", 2, 1, $false, $false, $false, $false) | Out-Null

# --- 2. Convert the used range into an Excel Table ---------------------------
$tableRange = $ws.Range("A1:U57")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row -------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
